# HW3/results_test1.xlsx -- "small but important fix"
# Re-derives the ensemble-size sweep (row 2 headers) and the accuracy
# numbers underneath it, clears the rows that no longer have data for
# this run, and restores the cursor/zoom to where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet view: zoom to 100% and move the selection -----------------
$ws.Select() | Out-Null
$excel.ActiveWindow.Zoom = 100
$ws.Range("G18").Select() | Out-Null

# --- row 2: ensemble-size header values shift down by one step -------
$row2 = @{
    "D2" = 1
    "E2" = 3
    "F2" = 5
    "G2" = 7
    "H2" = 9
    "I2" = 11
    "J2" = 13
    "K2" = 15
    "L2" = 17
    "M2" = 19
}
foreach ($addr in $row2.Keys) {
    $ws.Range($addr).Value = $row2[$addr]
}

# --- row 3 -------------------------------------------------------------
$row3 = @{
    "D3" = 0.75121951219512195
    "E3" = 0.81036585365853597
    "F3" = 0.83780487804878001
    "G3" = 0.85243902439024299
    "H3" = 0.85365853658536595
    "I3" = 0.85731707317073103
    "J3" = 0.85975609756097504
    "K3" = 0.86158536585365797
    "L3" = 0.86219512195121895
    "M3" = 0.86280487804878003
}
foreach ($addr in $row3.Keys) {
    $ws.Range($addr).Value = $row3[$addr]
}

# --- row 4 -------------------------------------------------------------
$row4 = @{
    "D4" = 0.86185574699783696
    "E4" = 0.86002088461251502
    "F4" = 0.86002088461251502
    "G4" = 0.86002088461251502
    "H4" = 0.86002088461251502
    "I4" = 0.86002088461251502
    "J4" = 0.86002088461251502
    "K4" = 0.86002088461251502
    "L4" = 0.86002088461251502
    "M4" = 0.86002088461251502
}
foreach ($addr in $row4.Keys) {
    $ws.Range($addr).Value = $row4[$addr]
}

# --- rows 5-8: no longer populated for this run -> clear ---------------
$ws.Range("D5:M8").ClearContents()

# --- row 9 ---------------------------------------------------------------
$row9 = @{
    "D9" = 0.64715909090908996
    "E9" = 0.803367582210749
    "F9" = 0.76006964050442305
    "G9" = 0.78319637143394905
    "H9" = 0.89612002903928301
    "I9" = 0.91108248635422495
    "J9" = 0.92060293484982902
    "K9" = 0.93049679357909099
    "L9" = 0.93235107284020302
    "M9" = 0.94088640396870205
}
foreach ($addr in $row9.Keys) {
    $ws.Range($addr).Value = $row9[$addr]
}

# --- row 10 (M10 stays blank, as in the source) --------------------------
$row10 = @{
    "D10" = 0.73003229948105697
    "E10" = 0.80744501357855303
    "F10" = 0.76019349439380501
    "G10" = 0.87187760479685905
    "H10" = 0.88993540103788504
    "I10" = 0.90996763329837804
    "J10" = 0.91615108493990405
    "K10" = 0.92184080153800596
    "L10" = 0.93717885375494003
}
foreach ($addr in $row10.Keys) {
    $ws.Range($addr).Value = $row10[$addr]
}
$ws.Range("M10").ClearContents()

# --- rows 11-14: no longer populated for this run -> clear ---------------
$ws.Range("D11:M14").ClearContents()
